{"js": "const replacements = [\n  [\"703\u00d76=4218\", \"442\u00d79=3978\"],\n  [\"333\u00d72=666\", \"271\u00d72=542\"],\n  [\"925\u00d74=3700\", \"769\u00d75=3845\"],\n  [\"958\u00d78=7664\", \"368\u00d77=2576\"],\n  [\"231\u00d76=1386\", \"743\u00d76=4458\"],\n  [\"530\u00d72=1060\", \"271\u00d75=1355\"],\n  [\"686\u00d74=2744\", \"573\u00d74=2292\"],\n  [\"353\u00d78=2824\", \"240\u00d78=1920\"],\n  [\"371\u00d74=1484\", \"326\u00d78=2608\"],\n  [\"239\u00d73=717\", \"343\u00d79=3087\"],\n  [\"830\u00d74=3320\", \"584\u00d78=4672\"],\n  [\"183\u00d72=366\", \"235\u00d79=2115\"],\n  [\"747\u00d78=5976\", \"487\u00d78=3896\"],\n  [\"601\u00d77=4207\", \"872\u00d72=1744\"],\n  [\"661\u00d79=5949\", \"181\u00d76=1086\"],\n  [\"837\u00d79=7533\", \"296\u00d79=2664\"],\n  [\"507\u00d79=4563\", \"226\u00d74=904\"],\n  [\"349\u00d78=2792\", \"913\u00d78=7304\"],\n  [\"773\u00d73=2319\", \"489\u00d72=978\"],\n  [\"172\u00d76=1032\", \"305\u00d74=1220\"],\n  [\"702\u00d72=1404\", \"736\u00d72=1472\"],\n  [\"172\u00d74=688\", \"271\u00d73=813\"],\n  [\"373\u00d72=746\", \"117\u00d76=702\"],\n  [\"380\u00d77=2660\", \"934\u00d73=2802\"],\n  [\"736\u00d76=4416\", \"704\u00d75=3520\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"703\u00d76=4218\", \"442\u00d79=3978\"),\n    @(\"333\u00d72=666\", \"271\u00d72=542\"),\n    @(\"925\u00d74=3700\", \"769\u00d75=3845\"),\n    @(\"958\u00d78=7664\", \"368\u00d77=2576\"),\n    @(\"231\u00d76=1386\", \"743\u00d76=4458\"),\n    @(\"530\u00d72=1060\", \"271\u00d75=1355\"),\n    @(\"686\u00d74=2744\", \"573\u00d74=2292\"),\n    @(\"353\u00d78=2824\", \"240\u00d78=1920\"),\n    @(\"371\u00d74=1484\", \"326\u00d78=2608\"),\n    @(\"239\u00d73=717\", \"343\u00d79=3087\"),\n    @(\"830\u00d74=3320\", \"584\u00d78=4672\"),\n    @(\"183\u00d72=366\", \"235\u00d79=2115\"),\n    @(\"747\u00d78=5976\", \"487\u00d78=3896\"),\n    @(\"601\u00d77=4207\", \"872\u00d72=1744\"),\n    @(\"661\u00d79=5949\", \"181\u00d76=1086\"),\n    @(\"837\u00d79=7533\", \"296\u00d79=2664\"),\n    @(\"507\u00d79=4563\", \"226\u00d74=904\"),\n    @(\"349\u00d78=2792\", \"913\u00d78=7304\"),\n    @(\"773\u00d73=2319\", \"489\u00d72=978\"),\n    @(\"172\u00d76=1032\", \"305\u00d74=1220\"),\n    @(\"702\u00d72=1404\", \"736\u00d72=1472\"),\n    @(\"172\u00d74=688\", \"271\u00d73=813\"),\n    @(\"373\u00d72=746\", \"117\u00d76=702\"),\n    @(\"380\u00d77=2660\", \"934\u00d73=2802\"),\n    @(\"736\u00d76=4416\", \"704\u00d75=3520\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    ) | Out-Null\n}"}
